$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.48%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.26%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.123"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.02%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07831"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.81%"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.250"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.53%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.891"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.92%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.965"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.87%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9217"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.33%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1081"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-9.75%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1894"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.06%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08868"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.58%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03347"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.75%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09570"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.24%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001377"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.58%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005668"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.44%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.413"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.79%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.397"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.06%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3424"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.70%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.281"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "19.59%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.51%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04361"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.76%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001194"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.52%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004265"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.34%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001399"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "16.44%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02178"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.88%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05033"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.06%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007533"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.09%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008665"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-11.83%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.13%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.007896"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-9.90%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006554"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.34%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.15%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "12.22%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-16.52%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.15%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.15%"
